# fixed test issue with name of list docs in delivery
# Adds a new "TestUserProd1" row to the Users sheet (row 28) mirroring the
# existing rows above it: UserName, Password, Email (as hyperlink) and
# Password/Answer columns, then updates the sheet's selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# --- Fill in the new user row (row 28) ---
$ws.Range("A28").Value = "TestUserProd1"
$ws.Range("B28").Value = "Password1"
$ws.Range("H28").Value = "thomsonreuters"

# Add the e-mail hyperlink for the new row (mirrors rows 26/27 above it).
$ws.Hyperlinks.Add($ws.Range("G28"), "mailto:TestUserProd1@mailinator.com")
$ws.Range("G28").Value = "TestUserProd1@mailinator.com"

# Re-apply the formatting used by the row above so the hyperlink cell keeps
# the same cell style as its neighbours instead of the default one the
# Hyperlinks.Add call produces.
$ws.Range("G27").Copy()
$ws.Range("G28").PasteSpecial(-4122)

# --- Update the active selection on the Users sheet ---
[void]$ws.Range("M10").Select()
